$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty scores for row 7 (DIONE / SAMBA)
$ws.Range("D7").Value = 8
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 7
$ws.Range("G7").Value = 7
$ws.Range("H7").Value = 7
$ws.Range("I7").Value = 8

# Correct the "Comportement" score for row 9 (DIOP / MAMADOU)
$ws.Range("H9").Value = 5

# Match the author's final selection/active cell
$ws.Range("H9").Select()
